# Update the "Estado de Cuenta" database: swap the two Periodo Mora rows
# (2505 / 2506) together with their corresponding "Valor Mora" amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (period 2506 / value 56940) becomes period 2505 / value 7592
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 7592

# Row 17 (period 2505 / value 7592) becomes period 2506 / value 56940
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
